$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "905×3=2715" "316×6=1896"
Replace-Text "883×2=1766" "348×3=1044"
Replace-Text "426×6=2556" "147×9=1323"
Replace-Text "365×3=1095" "572×2=1144"
Replace-Text "412×9=3708" "231×3=693"
Replace-Text "913×9=8217" "395×4=1580"
Replace-Text "290×7=2030" "773×9=6957"
Replace-Text "907×7=6349" "404×4=1616"
Replace-Text "699×9=6291" "324×2=648"
Replace-Text "863×8=6904" "339×9=3051"
Replace-Text "173×9=1557" "643×4=2572"
Replace-Text "269×9=2421" "599×9=5391"
Replace-Text "592×4=2368" "427×3=1281"
Replace-Text "263×4=1052" "660×5=3300"
Replace-Text "750×2=1500" "160×3=480"
Replace-Text "671×9=6039" "621×6=3726"
Replace-Text "421×6=2526" "730×6=4380"
Replace-Text "712×5=3560" "211×2=422"
Replace-Text "370×8=2960" "308×9=2772"
Replace-Text "995×6=5970" "807×9=7263"
Replace-Text "512×3=1536" "352×9=3168"
Replace-Text "330×3=990" "919×6=5514"
Replace-Text "847×8=6776" "550×7=3850"
Replace-Text "377×5=1885" "533×5=2665"
Replace-Text "170×9=1530" "823×6=4938"
